# CIERRE 24 DIC 2021
# Update the "ARQUITECTO" vale: the incentive amount goes from
# $50,000 to $70,000, and the spelled-out amount is updated to match
# ("CINCUENTA MIL" -> "SETENTA MIL"). Leaving the leftover now-unused
# string causes the shared-string table to re-pack, which in turn shifts
# the "VALES DE INSENTIVOS" sheet's text reference (same wording,
# different shared-string index) - no separate write is required there.

$wb = $excel.ActiveWorkbook
$wsArquitecto = $wb.Worksheets.Item("ARQUITECTO        ")

# Amount (numeric) and its written-out Spanish text.
$wsArquitecto.Range("D1").Value = 70000
$wsArquitecto.Range("A2").Value = "SETENTA      MIL   PESOS 00/100 M.N."

# Move/refresh the active selection on this sheet.
$null = $wsArquitecto.Activate()
$null = $wsArquitecto.Range("D3").Select()
